$wb = $excel.ActiveWorkbook

# --- 1) "hydro_catchment_data": append a new row 81 with the latest
#        (01/07/2025) reading, extracted from the climate update PDF and
#        already staged on the "Data" sheet (row 2). ---
$hydroSheet = $wb.Worksheets.Item("hydro_catchment_data")

$hydroSheet.Range("A81").NumberFormat = "@"
$hydroSheet.Range("A81").Value = "01/07/2025"
$hydroSheet.Range("B81:M81").Value = 0
$hydroSheet.Range("N81").Value = "(Rathnapura)"
$hydroSheet.Range("O81").Value = 0
$hydroSheet.Range("P81").Value = "NA"
$hydroSheet.Range("Q81").Value = "NA"

# --- 2) "Data": the staged row (row 2) now matches the header/formatted
#        row 1 above it, so line it up with the same cell formatting. ---
$dataSheet = $wb.Worksheets.Item("Data")

$dataSheet.Range("A2:Q2").Style = $dataSheet.Range("A1").Style

Write-Output "done"
